# Cập nhật data.xlsx từ công cụ QR
# The QR tool re-scanned and the duplicate/old entry in row 2 is removed,
# so the newer scan (previously row 3) shifts up to become row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("locations")

# Delete the old row 2 entirely; Excel shifts row 3 up into its place.
$ws.Rows.Item(2).Delete()
